$wb = $excel.ActiveWorkbook

# --- PRODUTOS sheet: rename table "Produto" -> "Produtos" in column L formulas,
#     and update the saved selection (was C38, now L2:L34) ---
$wsProdutos = $wb.Worksheets.Item("PRODUTOS")
$wsProdutos.Range("L2:L34").Replace("INSERT INTO Produto(", "INSERT INTO Produtos(", -4123)
$wsProdutos.Range("L2:L34").Select()

# --- FORNECEDOR sheet: rename table "Fornecedor" -> "Fornecedores" ---
$wsFornecedor = $wb.Worksheets.Item("FORNECEDOR")
$wsFornecedor.Range("E2:E5").Replace("INSERT INTO Fornecedor(", "INSERT INTO Fornecedores(", -4123)

# --- CATEGORIA sheet: rename table "Categoria" -> "Categorias" ---
$wsCategoria = $wb.Worksheets.Item("CATEGORIA")
$wsCategoria.Range("C2:C5").Replace("INSERT INTO Categoria(", "INSERT INTO Categorias(", -4123)

# --- SOLICITANTE sheet: rename table "Solicitante" -> "Solicitantes" ---
$wsSolicitante = $wb.Worksheets.Item("SOLICITANTE")
$wsSolicitante.Range("G2:G7").Replace("INSERT INTO Solicitante(", "INSERT INTO Solicitantes(", -4123)

# --- ENDERECO sheet: rename table "Endereco" -> "Enderecos" ---
$wsEndereco = $wb.Worksheets.Item("ENDERECO")
$wsEndereco.Range("H2:H11").Replace("INSERT INTO Endereco (", "INSERT INTO Enderecos (", -4123)

# --- Make CATEGORIA the active tab (this also clears tabSelected on PRODUTOS
#     and sets workbook activeTab) ---
$wsCategoria.Activate()
